$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "Förändrad" (changed) date in column C for every data row. ---
# Column C currently holds the Excel serial date 45184 (2023-09-15) and the
# commit bumps it to 45186 (2023-09-17) for every populated data row.
$lastRow = $ws.Cells.Item(1, 3).End(-4121).Row   # xlDown = -4121, walks from header (row1) to last used row in column C
if ($lastRow -lt 2) { $lastRow = 339 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# --- 2) Add a friendly display-text second argument to every HYPERLINK() formula. ---
# Every HYPERLINK(...) formula on the sheet (columns S,T,U,V,W,X,Y) gets the
# row's "Beteckning" (column A) appended as the link's friendly text, e.g.
#   HYPERLINK("...A 8503-2019.xlsx")  ->  HYPERLINK("...A 8503-2019.xlsx", "A 8503-2019")
$hyperlinkCols = 19, 20, 21, 22, 23, 24, 25   # S, T, U, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($beteckning)) { continue }

    foreach ($c in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $c)
        if (-not $cell.HasFormula) { continue }

        $f = $cell.Formula
        if ($f -notmatch '^=HYPERLINK\(') { continue }
        if ($f.TrimEnd() -match ',\s*"[^"]*"\)$') { continue }  # already has a friendly name

        $trimmed = $f.TrimEnd()
        $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $beteckning + '")'
        $cell.Formula = $newFormula
    }
}
